$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")

# Screenshot utility (flight-price scraper) test run logged its latest
# result into row 2 of the "Output" sheet: the completion timestamp plus
# the scraped fare breakdown.
$ws.Range("A2").Value = "21/01/2022 11:23:32 am"
$ws.Range("D2").Value = "₹1,21,740"
$ws.Range("E2").Value = "₹8,356"
$ws.Range("G2").Value = "₹1,30,106"
